$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the SUPPORTSTAFF table (columns D/E) with the new Janitor row ---
$ws.Range("D16").Style = "20% - Accent1"
$ws.Range("D16").Value = 7420
$ws.Range("E16").Style = "20% - Accent1"
$ws.Range("E16").Value = "Janitor"

# Annotation under the SUPPORTSTAFF table
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").Value = "2018 aussie janitor^^"

# --- Extend the TEAMMEMBER table (columns N/O) with the new member row ---
$ws.Range("N22").Style = "20% - Accent1"
$ws.Range("N22").Value = 7420
$ws.Range("O22").Style = "20% - Accent1"
$ws.Range("O22").Value = "Australian Janitor"

# Annotation under the TEAMMEMBER table
$ws.Range("N23").Style = "Normal"
$ws.Range("N23").Value = "2018 aussie janitor^^"

# --- Replace the old warning cell I38 with a new MATCH table row (I/J/K) ---
$ws.Range("I38").Style = "20% - Accent1"
$ws.Range("I38").Value = 2018
$ws.Range("J38").Style = "20% - Accent1"
$ws.Range("J38").Value = 78242
$ws.Range("K38").Style = "20% - Accent1"
$ws.Range("K38").Value = 7420

# Annotation under the MATCH table (new row)
$ws.Range("I39").Style = "Normal"
$ws.Range("I39").Value = "2018 aussie janitor^^"

# The "Bad" cell style is now unused (its only user, the old I38, was
# overwritten above) -- remove it like Excel does when it stops being
# referenced on save.
$wb.Styles.Item("Bad").Delete()

# --- Update the view: selection moved to the newly edited rows ---
$ws.Range("I38:K38").Select()
